$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "cert_list_test"

# Add the new data row (row 6)
$ws.Range("A6").Value = "Multiple Name"
$ws.Range("B6").Value = "Crazy Person"
$ws.Range("C6").Value = "mp@gmail.com"
$ws.Range("D6").Value = "Learning Python 101"
$ws.Range("E6").Value = 45234
$ws.Range("F6").Value = "Jeffry Python"

# Add hyperlink for the new email cell
[void]$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:mp@gmail.com")

# Re-apply the same formatting as the rows above (hyperlink style for C6,
# date style for D6:E6) so the new row reuses the existing cell styles
# instead of creating new ones
$ws.Range("C2").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("D2:E2").Copy()
$ws.Range("D6:E6").PasteSpecial(-4122)

# Update the active cell selection
[void]$ws.Range("D9").Select()
